$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 63 first (higher row number) so row indices below it aren't affected
# when we then delete row 18.
$ws.Rows.Item(63).Delete()
$ws.Rows.Item(18).Delete()
